# Apply the "add mgt fee ppt" edit to Sheet1:
#  - H10 previously held a (now incorrect) management-fee value; it is cleared.
#  - J10 previously held an empty-string placeholder; it now holds the numeric 0.
#  - The sheet's active selection moves from K8 to J11 (the next reviewed cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the stale value that used to live in H10.
$ws.Range("H10").ClearContents()

# J10 now carries an explicit numeric zero instead of the old empty shared string.
$ws.Range("J10").Value = 0

# Make sure Sheet1 is the active sheet, then move the selection to J11 to match
# the saved cursor position in the edited workbook.
$ws.Activate()
$ws.Range("J11").Select()
